# Scheduled data refresh: currentAveragePrice / Leve price & profit columns (H:N)
# across all job sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

# --- ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H31").Value = 4000
$ws.Range("I31").Value = 4000
$ws.Range("K31").Value = 12000
$ws.Range("M31").Value = -11770
$ws.Range("H33").Value = 151550.6
$ws.Range("I33").Value = 750750
$ws.Range("J33").Value = 1750.75
$ws.Range("K33").Value = 750750
$ws.Range("L33").Value = 1750.75
$ws.Range("M33").Value = -750521
$ws.Range("N33").Value = -2208.75
$ws.Range("H86").Value = 6094222
$ws.Range("I86").Value = 2789.0667
$ws.Range("J86").Value = 11170416
$ws.Range("K86").Value = 2789.0667
$ws.Range("L86").Value = 11170416
$ws.Range("M86").Value = -1666.0667
$ws.Range("N86").Value = -11172662
$ws.Range("H89").Value = 6094222
$ws.Range("I89").Value = 2789.0667
$ws.Range("J89").Value = 11170416
$ws.Range("K89").Value = 13945.3335
$ws.Range("L89").Value = 55852080
$ws.Range("M89").Value = -8329.333499999999
$ws.Range("N89").Value = -55863312
$ws.Range("H112").Value = 2062.3462
$ws.Range("J112").Value = 1905.5454
$ws.Range("L112").Value = 5716.6362
$ws.Range("N112").Value = -7932.6362

# --- ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 3256.8572
$ws.Range("I74").Value = 3249.875
$ws.Range("J74").Value = 3331.3333
$ws.Range("K74").Value = 3249.875
$ws.Range("L74").Value = 3331.3333
$ws.Range("M74").Value = -2375.875
$ws.Range("N74").Value = -5079.3333
$ws.Range("H77").Value = 3256.8572
$ws.Range("I77").Value = 3249.875
$ws.Range("J77").Value = 3331.3333
$ws.Range("K77").Value = 16249.375
$ws.Range("L77").Value = 16656.6665
$ws.Range("M77").Value = -11881.375
$ws.Range("N77").Value = -25392.6665
$ws.Range("H92").Value = 0
$ws.Range("J92").Value = 0
$ws.Range("L92").Value = 0
$ws.Range("N92").ClearContents()
$ws.Range("H132").Value = 1257.3914
$ws.Range("I132").Value = 1178.1818
$ws.Range("J132").Value = 3000
$ws.Range("K132").Value = 3534.5454
$ws.Range("L132").Value = 9000
$ws.Range("M132").Value = -1004.5454
$ws.Range("N132").Value = -14060

# --- BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 13781.096
$ws.Range("I20").Value = 17656
$ws.Range("K20").Value = 17656
$ws.Range("M20").Value = -17409
$ws.Range("H80").Value = 55555750
$ws.Range("I80").Value = 260.6
$ws.Range("K80").Value = 260.6
$ws.Range("M80").Value = 737.4
$ws.Range("H83").Value = 55555750
$ws.Range("I83").Value = 260.6
$ws.Range("K83").Value = 1303
$ws.Range("M83").Value = 3689

# --- CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 5422.149
$ws.Range("I31").Value = 11895.728
$ws.Range("J31").Value = 3444.111
$ws.Range("K31").Value = 11895.728
$ws.Range("L31").Value = 3444.111
$ws.Range("M31").Value = -11600.728
$ws.Range("N31").Value = -4034.111
$ws.Range("H34").Value = 5422.149
$ws.Range("I34").Value = 11895.728
$ws.Range("J34").Value = 3444.111
$ws.Range("K34").Value = 11895.728
$ws.Range("L34").Value = 3444.111
$ws.Range("M34").Value = -11693.728
$ws.Range("N34").Value = -3848.111
$ws.Range("H122").Value = 4390.467
$ws.Range("I122").Value = 4296.923
$ws.Range("J122").Value = 4998.5
$ws.Range("K122").Value = 12890.769
$ws.Range("L122").Value = 14995.5
$ws.Range("M122").Value = -10440.769
$ws.Range("N122").Value = -19895.5
$ws.Range("H141").Value = 94235.69
$ws.Range("I141").Value = 35148
$ws.Range("J141").Value = 102676.79
$ws.Range("K141").Value = 35148
$ws.Range("L141").Value = 102676.79
$ws.Range("M141").Value = -29968
$ws.Range("N141").Value = -113036.79

# --- CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 862916.7
$ws.Range("I4").Value = 1341798.9
$ws.Range("K4").Value = 4025396.7
$ws.Range("M4").Value = -4025284.7
$ws.Range("H5").Value = 986.04346
$ws.Range("I5").Value = 586.2308
$ws.Range("K5").Value = 1758.6924
$ws.Range("M5").Value = -1646.6924
$ws.Range("H48").Value = 500
$ws.Range("I48").Value = 500
$ws.Range("K48").Value = 1500
$ws.Range("M48").Value = -1250
$ws.Range("H113").Value = 1463.2106
$ws.Range("I113").Value = 780.1429000000001
$ws.Range("K113").Value = 2340.4287
$ws.Range("M113").Value = -170.4287000000004
$ws.Range("H135").Value = 986.04346
$ws.Range("I135").Value = 586.2308
$ws.Range("K135").Value = 5276.077200000001
$ws.Range("M135").Value = -2741.077200000001

# --- GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 2754
$ws.Range("I113").Value = 1890.1538
$ws.Range("K113").Value = 1890.1538
$ws.Range("M113").Value = 279.8462
$ws.Range("H132").Value = 25950
$ws.Range("I132").Value = 41956.637
$ws.Range("K132").Value = 125869.911
$ws.Range("M132").Value = -123339.911

# --- LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3060.1667
$ws.Range("I40").Value = 2423.75
$ws.Range("J40").Value = 4333
$ws.Range("K40").Value = 2423.75
$ws.Range("L40").Value = 4333
$ws.Range("M40").Value = -2287.75
$ws.Range("N40").Value = -4605
$ws.Range("H132").Value = 6531.227
$ws.Range("I132").Value = 3546.6
$ws.Range("K132").Value = 10639.8
$ws.Range("M132").Value = -8109.799999999999
$ws.Range("H136").Value = 2370
$ws.Range("I136").Value = 1802.9524
$ws.Range("K136").Value = 5408.857199999999
$ws.Range("M136").Value = -2858.857199999999

# --- WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H38").Value = 5000
$ws.Range("I38").Value = 5000
$ws.Range("K38").Value = 5000
$ws.Range("M38").Value = -4527
$ws.Range("H49").Value = 50000
$ws.Range("I49").Value = 50000
$ws.Range("K49").Value = 50000
$ws.Range("M49").Value = -49770
$ws.Range("H62").Value = 1345751
$ws.Range("I62").Value = 3974254.2
$ws.Range("J62").Value = 31499.25
$ws.Range("K62").Value = 3974254.2
$ws.Range("L62").Value = 31499.25
$ws.Range("M62").Value = -3973630.2
$ws.Range("N62").Value = -32747.25
$ws.Range("H65").Value = 1345751
$ws.Range("I65").Value = 3974254.2
$ws.Range("J65").Value = 31499.25
$ws.Range("K65").Value = 19871271
$ws.Range("L65").Value = 157496.25
$ws.Range("M65").Value = -19868151
$ws.Range("N65").Value = -163736.25
$ws.Range("H107").Value = 575.3333
$ws.Range("I107").Value = 573.25
$ws.Range("K107").Value = 1719.75
$ws.Range("M107").Value = 200.25
$ws.Range("H126").Value = 2035.6666
$ws.Range("I126").Value = 1992.8
$ws.Range("K126").Value = 5978.4
$ws.Range("M126").Value = -3508.4
$ws.Range("H132").Value = 3902.8647
$ws.Range("I132").Value = 4217.8335
$ws.Range("K132").Value = 12653.5005
$ws.Range("M132").Value = -10123.5005
